$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Account by year")

# Row 2 - Checking
$ws.Range("D2").Value = -116130.41
$ws.Range("E2").Value = -109552.69
$ws.Range("F2").Value = -122146.94
$ws.Range("G2").Value = -137516.7
$ws.Range("H2").Value = -32813.36
$ws.Range("I2").Value = -518160.1

# Row 3 - Savings
$ws.Range("D3").Value = 462151.42
$ws.Range("E3").Value = 449620.26
$ws.Range("F3").Value = 528792.92
$ws.Range("G3").Value = 409323.45
$ws.Range("H3").Value = 459273.54
$ws.Range("I3").Value = 2309161.59

# Row 4 - Mastercard (near-zero rounding artifacts with sign changes)
$ws.Range("D4").Value = -0.0
$ws.Range("E4").Value = 0.0
$ws.Range("F4").Value = 0.0
$ws.Range("G4").Value = -0.0
$ws.Range("H4").Value = -0.0
$ws.Range("I4").Value = -0.0

# Row 5 - Visa (near-zero rounding artifacts with sign changes)
$ws.Range("D5").Value = 0.0
$ws.Range("E5").Value = -0.0
$ws.Range("F5").Value = 0.0
$ws.Range("G5").Value = -0.0
$ws.Range("H5").Value = -0.0
$ws.Range("I5").Value = -0.0

# Row 6 - Total
$ws.Range("D6").Value = 346021.01
$ws.Range("E6").Value = 340067.57
$ws.Range("F6").Value = 406645.98
$ws.Range("G6").Value = 271806.75
$ws.Range("H6").Value = 426460.18
$ws.Range("I6").Value = 1791001.49
